# Weekly update: insert two new rows (for the most recent week) at the top
# of the Betarraga data block, pushing the existing rows down by two
# positions. The two newly inserted rows reuse the values of the prior
# week's "Primera"/"Segunda" pair (the two rows immediately above the
# insertion point), advancing the date by 7 days (one week).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceRow1 = 78
$sourceRow2 = 79
$insertAt   = 80

$lastCol = 18  # column R

# Capture the values of the two source rows (A:R) before we mutate anything.
$row1Values = @()
$row2Values = @()
for ($c = 1; $c -le $lastCol; $c++) {
    $row1Values += , $ws.Cells.Item($sourceRow1, $c).Value2
    $row2Values += , $ws.Cells.Item($sourceRow2, $c).Value2
}

# Push rows 80..204 down by two rows, creating two blank rows at 80:81.
$ws.Rows("$($insertAt):$($insertAt + 1)").Insert()

# Fill the two newly inserted rows with the captured values.
for ($c = 1; $c -le $lastCol; $c++) {
    $ws.Cells.Item($insertAt, $c).Value2 = $row1Values[$c - 1]
    $ws.Cells.Item($insertAt + 1, $c).Value2 = $row2Values[$c - 1]
}

# Advance the "Fecha" (column D = 4) of the two new rows by one week.
$ws.Cells.Item($insertAt, 4).Value2 = $row1Values[3] + 7
$ws.Cells.Item($insertAt + 1, 4).Value2 = $row2Values[3] + 7
